# fix status on sample sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "approved"
$ws.Range("D3").Value = "pending"
$ws.Range("D4").Value = "approved"

$ws.Range("E7").Select()
